$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values (written in this exact order so the shared-string table
#     ends up in the same sequence the target workbook uses) ---
$ws.Range("B2").Value = 'There are need to be "damage" instead of health and it''s needed to substract health with damage '
$ws.Range("C2").Value = "{ damage = amountOfDamage                                              health = health - damage;}"
$ws.Range("A3").Value = "public void resetDamage() { damage = damage }"
$ws.Range("A2").Value = "public void takeDamage(int amountOfDamage) { health = amountOfDamage }"
$ws.Range("B3").Value = "There are need to add line of code before    { damage = damage }"
$ws.Range("C3").Value = "health = health + getDamage();"

# --- Header row (A1:C1): make B1 match the same look as A1/C1 ---
$header = $ws.Range("A1:C1")
$header.Font.Bold = $true
$header.Font.Name = "Arial"
$header.Font.Size = 10
$header.Font.Color = 16777215
$header.Interior.Color = 6711008
$header.HorizontalAlignment = -4108
$header.WrapText = $true
$header.NumberFormat = "@"

# --- Data rows (A2:C3): center aligned, vertically centered text ---
foreach ($addr in @("A2","A3","C3")) {
    $c = $ws.Range($addr)
    $c.Font.Name = "Arial"
    $c.Font.Size = 10
    $c.Font.Color = 0
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4108
    $c.WrapText = $false
}
foreach ($addr in @("B2","C2","B3")) {
    $c = $ws.Range($addr)
    $c.Font.Name = "Arial"
    $c.Font.Size = 10
    $c.Font.Color = 0
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4108
    $c.WrapText = $true
}

# --- Row heights ---
$ws.Rows.Item(1).RowHeight = 15.75
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 30

# --- Column widths (closest achievable values for 60 / 38.1640625 / 41.6640625 chars) ---
$ws.Columns.Item(1).ColumnWidth = 59.166666666666664
$ws.Columns.Item(2).ColumnWidth = 37.333333333333336
$ws.Columns.Item(3).ColumnWidth = 40.833333333333336

# --- Remove the drawing/picture object left on the sheet ---
foreach ($shp in @($ws.Shapes)) {
    $shp.Delete()
}

# --- View: frozen header row, final selection C3, zoom 183% ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("C3").Select()
$excel.ActiveWindow.Zoom = 183
